$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: D3 no longer holds a numeric value (cell becomes blank)
$ws.Range("D3").ClearContents()

# Row 4: C4 corrected to 0
$ws.Range("C4").Value = 0

# Row 5: C5 corrected to 0
$ws.Range("C5").Value = 0

# Row 7: label renamed from "Other" to "Biogas", value corrected
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 35.90807867917246

# New row 8: re-add an "Other" row below Biogas, copying the label
# formatting used by the other row headers (bold, centered, bordered)
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = 34.64917958923758
